$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume snapshot (GitHub Actions bot update).
# Price ("D") cells contain text like "27.631.38" / "1.000" that Excel would
# otherwise silently reinterpret as a number (dropping the trailing zero /
# choking on the second "."), so those values are written with a leading
# apostrophe to force literal-text storage, exactly like the source sheet.
# Volume ("E") cells keep their "  +x.xx%  " padded text as-is.

$ws.Range("D2").Value = '''27.631.38'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '''1.843.20'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '''314.55'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '''0.4252'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").Value = '''0.3645'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = '''45.63'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = '''0.07275'
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("D11").Value = '''0.8986'
$ws.Range("E11").Value = '  -4.39%  '
$ws.Range("D12").Value = '''20.63'
$ws.Range("E12").Value = '  -3.83%  '
$ws.Range("D13").Value = '''1.818.65'
$ws.Range("E13").Value = '  -4.57%  '
$ws.Range("D14").Value = '''5.388'
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").Value = '''6.563'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").Value = '''0.06849'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").Value = '''78.33'
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").Value = '''0.000008869'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").Value = '''0.9997'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '''15.60'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").Value = '''27.630.74'
$ws.Range("E22").Value = '  -2.22%  '
$ws.Range("D23").Value = '''4.973'
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("E24").Value = '  -1.68%  '

# Ranks 23/24 swapped places on this refresh: Toncoin moved up to rank 23
# (row 25) and WrappedliquidstakedEther2.0 dropped to rank 24 (row 26). The
# "A" rank-index column itself is unchanged; only coin/link/price/volume move.
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''2.044'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '''1.988.99'
$ws.Range("E26").Value = '  -6.41%  '
$ws.Range("D27").Value = '''154.25'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").Value = '''18.26'
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").Value = '''5.250'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").Value = '''1.835'
$ws.Range("E30").Value = '  +6.10%  '
$ws.Range("D31").Value = '''110.82'
$ws.Range("E31").Value = '  -2.67%  '
$ws.Range("D32").Value = '''0.08871'
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").Value = '''0.7761'
$ws.Range("E33").Value = '  -2.88%  '
$ws.Range("D34").Value = '''4.567'
$ws.Range("E34").Value = '  -5.55%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Value = '''1.096'
$ws.Range("E36").Value = '  -6.45%  '
$ws.Range("D37").Value = '''0.9990'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '''1.097'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("D40").Value = '''0.01925'
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").Value = '''2.779'
$ws.Range("E41").Value = '  -6.16%  '
$ws.Range("D42").Value = '''0.5068'
$ws.Range("E42").Value = '  -3.55%  '
$ws.Range("D43").Value = '''6.797'
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("D44").Value = '''0.1641'
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").Value = '''8.227'
$ws.Range("E45").Value = '  -5.61%  '
$ws.Range("D46").Value = '''0.06636'
$ws.Range("E46").Value = '  -1.91%  '
$ws.Range("D47").Value = '''10.39'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").Value = '''0.4723'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = '''105.82'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").Value = '''0.9997'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '''1.637'
$ws.Range("E51").Value = '  -2.47%  '
